# Insert a new "SourceCode" paragraph (pander(table_forecasts) call) right
# after the "Using historical data, ... summarized in the table below:"
# paragraph and before the forecasts table that follows it.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "... summarized in the table below:"
$rng = $d.Content
$found = $rng.Find.Execute("table below:", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor text 'table below:' was not found in the document."
}

$anchorParagraph = $rng.Paragraphs(1)
$anchorRange = $anchorParagraph.Range

# Collapse to the end of that paragraph, then add a brand new paragraph
# right after it (this keeps it at body level, ahead of the <w:tbl> that
# immediately follows).
$anchorRange.Collapse(0)   # wdCollapseEnd
$anchorRange.InsertParagraphAfter()

# Grab a handle to the paragraph we just created.
$newParagraph = $anchorParagraph.Next()
$newParagraph.Style = "SourceCode"

$newRange = $newParagraph.Range
$newRange.Text = "pander(table_forecasts)"

$start = $newRange.Start

# First run: "pander" styled as FunctionTok
$functionRange = $d.Range($start, $start + 6)
$functionRange.Style = "FunctionTok"

# Second run: "(table_forecasts)" styled as NormalTok
$normalRange = $d.Range($start + 6, $start + 24)
$normalRange.Style = "NormalTok"

Write-Host "Inserted SourceCode paragraph: [$($newParagraph.Range.Text)]"
